$d = $word.ActiveDocument

# Locate the M2Doc field (a Word field written as " m:...' " instrText) that
# needs to be rewritten as literal "{m:...}" text runs instead of a real
# Word field (TokenIteratorFieldRewriterSplit no longer relies on Word
# fields to carry the template tokens).
$targetField = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields.Item($i)
    if ($candidate.Code.Text -like "*asTable*") {
        $targetField = $candidate
        break
    }
}

if ($targetField -eq $null) {
    throw "Could not find the asTable field to rewrite"
}

# Remember where the field lives so we can drop the literal text back in
# the same spot once the field (fldChar begin / instrText* / fldChar end)
# is gone. The field's owning paragraph starts right where the "begin"
# fldChar run used to sit, which is exactly where the replacement text
# should be inserted.
$fieldCodeStart = $targetField.Code.Start
$insertAt = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidatePara = $d.Paragraphs.Item($i)
    if (($candidatePara.Range.Start -le $fieldCodeStart) -and ($fieldCodeStart -le $candidatePara.Range.End)) {
        $insertAt = $candidatePara.Range.Start
        break
    }
}

if ($insertAt -eq $null) {
    throw "Could not locate the paragraph owning the asTable field"
}

# Deleting the field removes the begin/end fldChar runs together with all
# of the instrText runs, leaving an empty paragraph (pPr is preserved).
$targetField.Delete()

# Rebuild the same textual content as plain "w:t" runs - one run per
# former instrText run, each keeping its xml:space="preserve" flag - but
# wrapped as literal "{m:...}" text instead of a hidden field code.
$runsXml = ''
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{m:</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + "'excel.xlsx'.asTable('" + '</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Feuil1</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">' + "', " + '</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + "'C3', 'F7'" + '</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + ", 'fr" + '</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>' + "FR'" + '</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>'
$runsXml += '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">}</w:t></w:r>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
    $runsXml +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Range($insertAt, $insertAt)
$target.InsertXML($packageXml)

Write-Output ("Paragraph now reads: " + $d.Paragraphs.Item(2).Range.Text)
